# Completed end to end scenario of product checkout:
# adds payee/card-payment columns (S:W) to the "EndToEnd" sheet and
# swaps the selected product/price in the existing cart columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("EndToEnd")

# --- Update the already-selected product (cart) to "Summer White Top" / Rs. 400 ---
$ws.Range("K2").Value = "Summer White Top"
$ws.Range("L2").Value = "Rs. 400"
$ws.Range("N2").Value = "Rs. 400"
$ws.Range("R2").Value = "Rs. 900"

# --- Add new payment-detail header cells (S1:W1), matching the header style used by R1 ---
$ws.Range("R1").Copy()
$ws.Range("S1:W1").PasteSpecial(-4122)

$ws.Range("S1").Value = "payeename"
$ws.Range("T1").Value = "cardnum"
$ws.Range("U1").Value = "cvc"
$ws.Range("V1").Value = "monthexpiry"
$ws.Range("W1").Value = "yearexpiry"

# --- Add the corresponding payment-detail data cells (S2:W2) ---
$ws.Range("S2").Value = "Tester"
$ws.Range("T2").Value = 12345
$ws.Range("U2").Value = 311
$ws.Range("V2").Value = 11
$ws.Range("W2").Value = 2024

# --- Update the sheet view: scroll position and the active selection ---
$ws.Activate()
$win = $excel.ActiveWindow
$win.ScrollColumn = 11
$win.ScrollRow = 2
$ws.Range("V7").Select()
